$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2..19 down to 3..20)
$ws.Rows("2:2").Insert()

# Set the new starting intensity value
$ws.Range("A2").Value = 1

# Update selection to match the target state (activeCell B7)
$ws.Range("B7").Select()
